$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 129 (existing rows 129..218 shift down to 131..220)
$ws.Rows("129:130").Insert()

# New row 129: Pomelo, Start Ruby, Primera, new weekly price record
$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44596
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100102
$ws.Range("H129").Value = "Cítricos"
$ws.Range("I129").Value = 100102006
$ws.Range("J129").Value = "Pomelo"
$ws.Range("K129").Value = "Start Ruby"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 180
$ws.Range("N129").Value = 13000
$ws.Range("O129").Value = 14000
$ws.Range("P129").Value = 13500
$ws.Range("Q129").Value = "$/caja 14 kilos empedrada"
$ws.Range("R129").Value = "Región de O'Higgins"
$ws.Range("S129").Value = 964
$ws.Range("T129").Value = 14

# New row 130: Pomelo, Start Ruby, Segunda, new weekly price record (same date as row 129)
$ws.Range("A130").Value = 4
$ws.Range("B130").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C130").Value = "Los Lagos"
$ws.Range("D130").Value = 44596
$ws.Range("E130").Value = 10
$ws.Range("F130").Value = "Fruta"
$ws.Range("G130").Value = 100102
$ws.Range("H130").Value = "Cítricos"
$ws.Range("I130").Value = 100102006
$ws.Range("J130").Value = "Pomelo"
$ws.Range("K130").Value = "Start Ruby"
$ws.Range("L130").Value = "Segunda"
$ws.Range("M130").Value = 60
$ws.Range("N130").Value = 11000
$ws.Range("O130").Value = 11000
$ws.Range("P130").Value = 11000
$ws.Range("Q130").Value = "$/caja 14 kilos empedrada"
$ws.Range("R130").Value = "Región de O'Higgins"
$ws.Range("S130").Value = 786
$ws.Range("T130").Value = 14
